$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new "Identifier" column is inserted as the first column in the three
# "header + data" mini-tables of this sheet (rows 6-7, 12-13 and 18-19).
# The existing Code/Description/Space (etc.) columns shift one place to the
# right (B->C, C->D) and a bold "Identifier" header plus the corresponding
# full settings path is written into the freed-up column A.
# ---------------------------------------------------------------------------

function Shift-RowFormats($row) {
    # Push existing cell formatting one column to the right so the moved
    # values keep looking the way they used to (rightmost column first so
    # nothing gets clobbered before it is copied).
    $ws.Range("C$row").Copy()
    $ws.Range("D$row").PasteSpecial(-4122)
    $ws.Range("B$row").Copy()
    $ws.Range("C$row").PasteSpecial(-4122)
    $ws.Range("A$row").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
}

function Set-IdentifierHeader($row) {
    $ws.Cells.Item($row, 1).Value = "Identifier"
    # New bold style (matches the other bold headers, but without the
    # left-alignment carried over from column A's previous formatting).
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 1).Font.Bold = $true
    $ws.Cells.Item($row, 1).Font.Size = 12
}

# --- PROJECT block (rows 6-7) ----------------------------------------------
Shift-RowFormats 6
Shift-RowFormats 7

$ws.Cells.Item(6,4).Value = "Space"
$ws.Cells.Item(6,3).Value = "Description"
$ws.Cells.Item(6,2).Value = "Code"
Set-IdentifierHeader 6

$ws.Cells.Item(7,4).Value = "ELN_SETTINGS"
$ws.Cells.Item(7,3).Value = "Default Project Updated"
$ws.Cells.Item(7,2).Value = "DEFAULT_PROJECT"
$ws.Cells.Item(7,1).Value = "/ELN_SETTINGS/DEFAULT_PROJECT"

# --- EXPERIMENT block (rows 12-13) ------------------------------------------
Shift-RowFormats 12
Shift-RowFormats 13

$ws.Cells.Item(12,4).Value = "Name"
$ws.Cells.Item(12,3).Value = "Project"
$ws.Cells.Item(12,2).Value = "Code"
Set-IdentifierHeader 12

$ws.Cells.Item(13,4).Value = "Default Experiment Updated"
$ws.Cells.Item(13,3).Value = "/ELN_SETTINGS/DEFAULT_PROJECT"
$ws.Cells.Item(13,2).Value = "DEFAULT_EXPERIMENT"
$ws.Cells.Item(13,1).Value = "/ELN_SETTINGS/DEFAULT_PROJECT/DEFAULT_EXPERIMENT"

# --- SAMPLE block (rows 18-19) ----------------------------------------------
Shift-RowFormats 18
Shift-RowFormats 19

$ws.Cells.Item(18,4).Value = "ELN Settings"
$ws.Cells.Item(18,3).Value = "Space"
$ws.Cells.Item(18,2).Value = "Code"
Set-IdentifierHeader 18

$ws.Cells.Item(19,4).Value = "{}"
$ws.Cells.Item(19,3).Value = "ELN_SETTINGS"
$ws.Cells.Item(19,2).Value = "GENERAL_ELN_SETTINGS"
$ws.Cells.Item(19,1).Value = "/ELN_SETTINGS/GENERAL_ELN_SETTINGS"

# --- Column widths (A widened for the long identifiers, B bumped up) -------
$ws.Columns.Item(1).ColumnWidth = 51.166666666666664
$ws.Columns.Item(2).ColumnWidth = 30.5

# --- Selection moves to A19 --------------------------------------------------
$ws.Range("A19").Select()
